$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Wenliang%Song%songwl2002jason@163.com%1,                    Junhua%Li%m13797631535@163.com%1,                    Ning%Zou%xdzning0215@163.com%1,                    Wenhe%Guan%guanwenhe@126.com%1,                    Jiali%Pan%cmu107819@163.com%1,                    Wei%Xu%tomxu.123@163.com%1]'
$ws.Range("I2").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E3").Value = '[Li%Zhu%NULL%1,                    Jian%Wang%NULL%0,                    Rui%Huang%NULL%1,                    Longgen%Liu%NULL%1,                    Haiyan%Zhao%NULL%1,                    Chao%Wu%dr.wu@nju.edu.cn%0,                    Chuanwu%Zhu%zhuchw@126.com%2,                    Chuanwu%Zhu%zhuchw@126.com%0]'
$ws.Range("I3").Value = '_PMC_CROSSREF'

$ws.Range("E4").Value = '[Mengqi%Liu%NULL%1,                    Zongbiao%Song%NULL%1,                    Kaihu%Xiao%NULL%0]'
$ws.Range("I4").Value = '_PMC_CROSSREF'

$ws.Range("E5").Value = '[Yang%Li%NULL%0,                    Jianghui%Cao%NULL%1,                    Xiaolong%Zhang%NULL%1,                    Guangzhi%Liu%NULL%1,                    Xiaxia%Wu%NULL%1,                    Baolin%Wu%WBLlin0129@163.com%1]'
$ws.Range("I5").Value = '_PMC_Springer_CROSSREF'

$ws.Range("E6").Value = '[Sharon%Steinberger%xref no email%1,            Bin%Lin%xref no email%1,            Adam%Bernheim%xref no email%1,            Michael%Chung%xref no email%1,            Yuantong%Gao%xref no email%1,            Zongyu%Xie%xref no email%1,            Tongtong%Zhao%xref no email%1,            Junli%Xia%xref no email%1,            Xueyan%Mei%xref no email%1,            Brent P.%Little%xref no email%1]'
$ws.Range("I6").Value = '_CROSSREF'

$ws.Range("E7").Value = '[Y.%Lu%NULL%1,                    H.%Wen%NULL%1,                    D.%Rong%NULL%1,                    Z.%Zhou%NULL%1,                    H.%Liu%NULL%1]'
$ws.Range("I7").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E8").Value = '[Lan%Lan%NULL%1,                    Dan%Xu%NULL%0,                    Dan%Xu%NULL%0,                    Chen%Xia%NULL%2,                    Chen%Xia%NULL%0,                    Shaokang%Wang%NULL%2,                    Shaokang%Wang%NULL%0,                    Minhua%Yu%NULL%2,                    Minhua%Yu%NULL%0,                    Haibo%Xu%NULL%0,                    Haibo%Xu%NULL%0]'
$ws.Range("I8").Value = '_PMC_CROSSREF'

$ws.Range("E9").Value = '[Jafar%Soltani%NULL%1,                    Iraj%Sedighi%NULL%1,                    Zohreh%Shalchi%NULL%1,                    Ghazal%Sami%NULL%1,                    Borhan%Moradveisi%NULL%1,                    Soheila%Nahidi%NULL%1]'
$ws.Range("I9").Value = '_PMC_CROSSREF'
